# Update the "dSF" column (F) values for specific rows on Sheet1.
# These edits reflect a repull/recalculation of data ("repull data, push all
# data, mean calculation") where only the dSF column values changed for a
# subset of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    13 = 0
    15 = -3
    16 = -1
    18 = 1
    24 = 0
    29 = 0
    38 = 4
    41 = -3
    46 = -1
    47 = -7
    54 = -4
    55 = -3
    56 = 5
    62 = 6
    64 = -7
    67 = -4
    69 = 3
    71 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
